# add-steps-in-new-test:
# Append six new "addStep" rows (9-14) to the Case sheet, matching the
# new UI-actions recorded for suite-1 / test-1, then widen column D to
# fit the longer parameter text and leave the selection on the last
# cell that was filled in (D14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Test", "test-1", "addStep", "object=CasesTab;name=casesTab;action=addSuite;arguments=suiteName:suite-1"),
    @("Test", "test-1", "addStep", "object=CasesTab;name=casesTab;action=openSuite;arguments=suiteName:suite-1"),
    @("Test", "test-1", "addStep", "object=Suite;name=suite-1;action=addTest;arguments=testName:test-1"),
    @("Test", "test-1", "addStep", "object=Suite;name=suite-1;action=openTest;arguments=testName:test-1"),
    @("Test", "test-1", "addStep", "object=Test;name=test-1;action=hasTotalSteps;arguments=numberOfSteps:0"),
    @("Test", "test-1", "addStep", "object=Test;name=test-1;action=addStep;arguments=object:Tool, name:tool, action:toCasesTab")
)

$row = 9
foreach ($rowValues in $newRows) {
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
    $ws.Cells.Item($row, 4).Value = $rowValues[3]
    $row = $row + 1
}

# Row 9 picked up a slightly taller custom row height when the data was
# imported.
$ws.Rows.Item(9).RowHeight = 14.25

# Column D needs to be widened so the long "addStep" parameter strings
# are fully visible.
$ws.Columns.Item(4).ColumnWidth = 88.85546875

# Leave the selection where the import left off.
$ws.Range("D14").Select() | Out-Null
